$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows before row 256 (shifts old rows 256-275 down to 258-277)
$ws.Rows.Item(256).Resize(2).Insert()

# --- New row 256: Primera, 29-jun-2023 (serial 45106) ---
$ws.Cells.Item(256, 1).Value = 7
$ws.Cells.Item(256, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(256, 3).Value = "Ñuble"
$ws.Cells.Item(256, 4).Value = 45106
$ws.Cells.Item(256, 5).Value = 16
$ws.Cells.Item(256, 6).Value = 100112040
$ws.Cells.Item(256, 7).Value = "Cilantro"
$ws.Cells.Item(256, 8).Value = "Sin especificar"
$ws.Cells.Item(256, 9).Value = "Primera"
$ws.Cells.Item(256, 10).Value = 180
$ws.Cells.Item(256, 11).Value = 1500
$ws.Cells.Item(256, 12).Value = 1500
$ws.Cells.Item(256, 13).Value = 1500
$ws.Cells.Item(256, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(256, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(256, 16).Value = 1500
$ws.Cells.Item(256, 17).Value = 1
$ws.Cells.Item(256, 18).Value = "Hortaliza"

# --- New row 257: Segunda, 29-jun-2023 (serial 45106) ---
$ws.Cells.Item(257, 1).Value = 7
$ws.Cells.Item(257, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(257, 3).Value = "Ñuble"
$ws.Cells.Item(257, 4).Value = 45106
$ws.Cells.Item(257, 5).Value = 16
$ws.Cells.Item(257, 6).Value = 100112040
$ws.Cells.Item(257, 7).Value = "Cilantro"
$ws.Cells.Item(257, 8).Value = "Sin especificar"
$ws.Cells.Item(257, 9).Value = "Segunda"
$ws.Cells.Item(257, 10).Value = 220
$ws.Cells.Item(257, 11).Value = 1000
$ws.Cells.Item(257, 12).Value = 1000
$ws.Cells.Item(257, 13).Value = 1000
$ws.Cells.Item(257, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(257, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(257, 16).Value = 1000
$ws.Cells.Item(257, 17).Value = 1
$ws.Cells.Item(257, 18).Value = "Hortaliza"
